$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 289, shifting rows 289:325 down to 290:326.
$ws.Rows(289).Insert()

# Populate the new row 289 with a new price record.
$ws.Cells.Item(289, 1).Value = 8
$ws.Cells.Item(289, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(289, 3).Value = "Coquimbo"
$ws.Cells.Item(289, 4).Value = 45077
$ws.Cells.Item(289, 5).Value = 4
$ws.Cells.Item(289, 6).Value = 100112037
$ws.Cells.Item(289, 7).Value = "Cebollín"
$ws.Cells.Item(289, 8).Value = "Sin especificar"
$ws.Cells.Item(289, 9).Value = "Primera"
$ws.Cells.Item(289, 10).Value = 600
$ws.Cells.Item(289, 11).Value = 5000
$ws.Cells.Item(289, 12).Value = 5500
$ws.Cells.Item(289, 13).Value = 5250
$ws.Cells.Item(289, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(289, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(289, 16).Value = 146
$ws.Cells.Item(289, 17).Value = 36
$ws.Cells.Item(289, 18).Value = "Hortaliza"
